# Adds Denmark, Sweden and Norway market sheets (cloned from the existing
# Belgium sheet's layout/styles) to the workbook, fills in their
# market-specific values, and updates the active-sheet/selection state so
# Norway ends up as the active (last) tab.

$wb = $excel.ActiveWorkbook

# Use "Belgium" as the template: same column widths, styles, merged cells
# and row layout are reused for every new market sheet.
$template = $wb.Worksheets.Item("Belgium")

# --- Denmark -----------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $afterSheet)
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-3446/T2010"
$denmark.Range("A1:XFD1048576").Select() | Out-Null

# --- Sweden --------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Copy($null, $afterSheet)
$sweden = $wb.Worksheets.Item($wb.Worksheets.Count)
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").Value = "NGC-3465/T2028"
$sweden.Range("A1:XFD1048576").Select() | Out-Null

# --- Norway --------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sweden.Copy($null, $afterSheet)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-3464/T1933"
$norway.Range("B5").Select() | Out-Null
